# "Alter bei AN-Anteil PV hinzugefügt und notwendige Änderungen vorgenommen."
#
# Insert a new question row ("juenger als 23 oder vor 1940 geboren?" / "nein")
# directly above the existing "wohnhaft Sachsen?" row (row 41), pushing every
# row below it down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Shift rows 41:48 down to 42:49 and open up a blank row 41.
$ws.Rows.Item(41).Insert() | Out-Null

# Populate the new row. Column B mirrors the numeric-style cell (the same
# style used directly above it, by "Anzahl Kinder") while holding the text
# answer "nein".
$ws.Cells.Item(41, 1).Value = "juenger als 23 oder vor 1940 geboren?"
$ws.Cells.Item(41, 2).Value = "nein"
$ws.Cells.Item(41, 2).NumberFormat = "0"

# Restore the view: scrolled down a bit with A28 as the active selection.
$ws.Range("A28").Select() | Out-Null
